$d = $word.ActiveDocument

# Locate the run "Both unit tests and system tests" (italic run at the end
# of the "Testing product" bullet).
$target = $d.Content
$found = $target.Find.Execute("Both unit tests and system tests", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Grab a single already-italic(+iCs) space character living *inside* that
    # same run (the space between "Both" and "unit") so we can clone its
    # run-properties exactly (xml:i + xml:iCs) onto a brand new run without
    # triggering the engine's run-coalescing pass (which fires on in-place
    # text edits, not on fresh FormattedText insert).
    $italicSpaceSource = $d.Range($target.Start + 4, $target.Start + 5)

    # Collapse to the point right after "...system tests" - this is where
    # the new content gets inserted.
    $target.Collapse(0)
    $insertPoint = $target.End

    # 1) Insert the italic space run.
    $spaceDest = $d.Range($insertPoint, $insertPoint)
    $spaceDest.FormattedText = $italicSpaceSource.FormattedText

    # 2) Insert the red-colored run right after the space.
    $redStart = $insertPoint + 1
    $redText = "as well as testing based on what user need access (student or professor)"
    $redDest = $d.Range($redStart, $redStart)
    $redDest.InsertAfter($redText)

    $redRange = $d.Range($redStart, $redStart + $redText.Length)
    $redRange.Font.Color = 255
}
